$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

# Status (row 6) : "active" -> "draft"
$meta.Range("B6").Value = "draft"

# Date (row 8) : republish timestamp
$meta.Range("B8").Value = "2023-08-01T16:12:28+00:00"
